# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 06:51"

# 2. Row 5 (country total update) - values refreshed
$ws.Range("B5").Value = 6685082
$ws.Range("C5").Value = 3009
$ws.Range("D5").Value = 5662490
$ws.Range("E5").Value = 918992

# 3. Thailand ("Tailandia") overtakes Gambia in the sorted list (by Casos totales),
#    so row 141 now holds Tailandia's refreshed figures and row 142 now holds
#    Gambia's (previously row 141's) unchanged figures.
$ws.Range("A141").Value = "Tailandia"
$ws.Range("B141").Value = 3600
$ws.Range("C141").Value = 10
$ws.Range("D141").Value = 3390
$ws.Range("E141").Value = 151
$ws.Range("H141").Value = 59

$ws.Range("A142").Value = "Gambia"
$ws.Range("B142").Value = 3594
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 2226
$ws.Range("E142").Value = 1253
$ws.Range("H142").Value = 115

# 4. Butan (row 187) and Camboya (row 188) figures refreshed
$ws.Range("B187").Value = 299
$ws.Range("C187").Value = 1
$ws.Range("D187").Value = 248
$ws.Range("E187").Value = 51

$ws.Range("D188").Value = 276
$ws.Range("E188").Value = 4
